$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "49.516.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.548.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.01"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.49"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.527"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.553"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.07"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.34"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0816"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.26"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.940.66"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.537.71"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.859"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "49.327.53"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.04"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +10.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.19"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.66"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0944"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "284.21"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.92"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.52"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.35"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +6.04%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.81"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.11"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.49"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.62"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.38"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0784"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.01"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.68"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.99"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.43"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.21"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.07"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.25%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.34"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +6.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.013.44"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.00"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.13"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +7.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.99"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.30"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.15"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.69%  "
